# The presentation's applied design is switched from the custom "Integral"
# (Red Violet) theme to the built-in default "Office Theme" (Office colour
# scheme) -- i.e. the Design/Theme gallery selection for the whole deck
# changes from the red-violet palette to the plain Office palette.
#
# PowerPoint stores each theme's twelve colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) on the active design's ColorScheme, which is
# reachable from the slide master (and mirrors onto every slide/notes/
# handout master in this deck, since they all resolve to the one live
# design). Re-pointing every one of the twelve slots to the stock "Office"
# palette reproduces the colour swap.

function PackRGB($r, $g, $b) {
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

# Target palette: the stock "Office" colour scheme (what the built-in
# "Office Theme" design uses), replacing the previous "Red Violet" colours.
$colorScheme.Item(1).RGB  = PackRGB 0x00 0x00 0x00   # dk1
$colorScheme.Item(2).RGB  = PackRGB 0xFF 0xFF 0xFF   # lt1
$colorScheme.Item(3).RGB  = PackRGB 0x44 0x54 0x6A   # dk2
$colorScheme.Item(4).RGB  = PackRGB 0xE7 0xE6 0xE6   # lt2
$colorScheme.Item(5).RGB  = PackRGB 0x5B 0x9B 0xD5   # accent1
$colorScheme.Item(6).RGB  = PackRGB 0xED 0x7D 0x31   # accent2
$colorScheme.Item(7).RGB  = PackRGB 0xA5 0xA5 0xA5   # accent3
$colorScheme.Item(8).RGB  = PackRGB 0xFF 0xC0 0x00   # accent4
$colorScheme.Item(9).RGB  = PackRGB 0x44 0x72 0xC4   # accent5
$colorScheme.Item(10).RGB = PackRGB 0x70 0xAD 0x47   # accent6
$colorScheme.Item(11).RGB = PackRGB 0x05 0x63 0xC1   # hlink
$colorScheme.Item(12).RGB = PackRGB 0x95 0x4F 0x72   # folHlink

# Carry the design/theme naming over to match ("Integral"/"Red Violet" ->
# "Office Theme"/"Office"); harmless if the host treats these as read-only.
try { $p.Designs.Item(1).Name = "Office Theme" } catch {}
try { $master.Theme.Name = "Office Theme" } catch {}
try { $colorScheme.Name = "Office" } catch {}
